# The measurement-error model gained "D" (delta) variants of the four
# count columns (SUSJ/RECJ/SUSA/RECA -> DSUSJ/DRECJ/DSUSA/DRECA) to match
# the revised C dmeasure/rmeasure code. "cumulative_time" (column A)
# stays the same; only the B1:E1 header labels change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "DSUSJ"
$ws.Range("C1").Value = "DRECJ"
$ws.Range("D1").Value = "DSUSA"
$ws.Range("E1").Value = "DRECA"

# The longer header labels now wrap onto a second line, so the header
# row grows from 15pt to 30pt.
$ws.Rows.Item(1).RowHeight = 30

Write-Output "header relabeled to D-prefixed measurement-error columns"
